$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.056.19"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.01"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.21%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.34"
$ws.Range("D5").ClearFormats()

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6355"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.39%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.77"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.87%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2938"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07348"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.85"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07679"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.829.51"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.22%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.989"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6637"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.97"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.69%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008678"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.22%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.065"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "28.912.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.96%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.079.48"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.115"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.52%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.465"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.61%  "

# Row 28
$ws.Range("E28").Value = "  -1.12%  "

# Row 29
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("E30").Value = "  -0.37%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.097"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.027"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.202"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.98%  "

# Row 34
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.834"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7386"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.152"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.294.50"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.21%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01784"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.56%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.746"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.17%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.301"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.46%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8964"
$ws.Range("D43").ClearFormats()

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.978.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5137"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.01"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.51%  "

# Row 49
$ws.Range("E49").Value = "  -9.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05828"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.70%  "
